# ---------------------------------------------------------------------------
# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2310"
#   "<name>_new" -> "<name>_FV2404"
# and wrap the sheet's data range (A1:U66) in a native Excel Table, plus
# freeze the header row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the 20 header cells (column K / "diff" stays unchanged).
# ---------------------------------------------------------------------------
$oldHeaders = @(
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310"
)
$newHeaders = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# ---------------------------------------------------------------------------
# 2) Turn A1:U66 into a native table ("Table1") without disturbing the
#    existing header-row formatting (creating the ListObject directly on
#    top of the already-styled header row would otherwise synthesize an
#    extra differential style for the header). We build the table on a
#    scratch range first and then resize/move it onto the real data.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("W1:X2")
$ws.Range("W1").Value = "Scratch1"
$ws.Range("X1").Value = "Scratch2"
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 2

$listObjects = $ws.ListObjects()
$table = $listObjects.Add(1, $scratch, $null, 1)
$table.TableStyle = ""

$table.Resize($ws.Range("A1:U66"))

# Clean up the scratch cells again now that the table lives on A1:U66.
$ws.Range("W1:X2").ClearContents()

# Re-sync the table's column headers with the real header-row text (writing
# through the header cell keeps styling untouched and updates the table
# column name at the same time).
for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
$ws.Cells.Item(1, 11).Value = "diff"
for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# ---------------------------------------------------------------------------
# 3) Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit complete"
